# Auto-generated edit script applying numeric updates to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H-N)
# across all 8 worksheets, per the scheduled-runner refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1231.5714
$ws.Cells.Item(19, 9).Value = 844.75
$ws.Cells.Item(19, 10).Value = 1747.3334
$ws.Cells.Item(19, 11).Value = 844.75
$ws.Cells.Item(19, 12).Value = 1747.3334
$ws.Cells.Item(19, 13).Value = -669.75
$ws.Cells.Item(19, 14).Value = -2097.3334

$ws.Cells.Item(51, 8).Value = 8678.799999999999
$ws.Cells.Item(51, 10).Value = 8897.462
$ws.Cells.Item(51, 12).Value = 8897.462
$ws.Cells.Item(51, 14).Value = -9865.462

$ws.Cells.Item(88, 8).Value = 1469.5625
$ws.Cells.Item(88, 9).Value = 741.8
$ws.Cells.Item(88, 10).Value = 1800.3636
$ws.Cells.Item(88, 11).Value = 741.8
$ws.Cells.Item(88, 12).Value = 1800.3636
$ws.Cells.Item(88, 13).Value = -335.8
$ws.Cells.Item(88, 14).Value = -2612.3636

$ws.Cells.Item(91, 8).Value = 1469.5625
$ws.Cells.Item(91, 9).Value = 741.8
$ws.Cells.Item(91, 10).Value = 1800.3636
$ws.Cells.Item(91, 11).Value = 741.8
$ws.Cells.Item(91, 12).Value = 1800.3636
$ws.Cells.Item(91, 13).Value = 662.2
$ws.Cells.Item(91, 14).Value = -4608.3636

$ws.Cells.Item(132, 8).Value = 333436.06
$ws.Cells.Item(132, 9).Value = 416732.53
$ws.Cells.Item(132, 11).Value = 1250197.59
$ws.Cells.Item(132, 13).Value = -1247667.59

$ws.Cells.Item(137, 8).Value = 9444.526
$ws.Cells.Item(137, 9).Value = 4824.4165
$ws.Cells.Item(137, 11).Value = 14473.2495
$ws.Cells.Item(137, 13).Value = -11923.2495

$ws.Cells.Item(138, 8).Value = 2792.768
$ws.Cells.Item(138, 10).Value = 5018.303
$ws.Cells.Item(138, 12).Value = 15054.909
$ws.Cells.Item(138, 14).Value = -25334.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1821.8125
$ws.Cells.Item(45, 9).Value = 1580.7693
$ws.Cells.Item(45, 10).Value = 2866.3333
$ws.Cells.Item(45, 11).Value = 1580.7693
$ws.Cells.Item(45, 12).Value = 2866.3333
$ws.Cells.Item(45, 13).Value = -1203.7693
$ws.Cells.Item(45, 14).Value = -3620.3333

$ws.Cells.Item(110, 8).Value = 6021.5654
$ws.Cells.Item(110, 9).Value = 2981.6667
$ws.Cells.Item(110, 10).Value = 7094.4707
$ws.Cells.Item(110, 11).Value = 2981.6667
$ws.Cells.Item(110, 12).Value = 7094.4707
$ws.Cells.Item(110, 13).Value = -936.6667000000002
$ws.Cells.Item(110, 14).Value = -11184.4707

$ws.Cells.Item(132, 8).Value = 1083889.5
$ws.Cells.Item(132, 9).Value = 1379389.5
$ws.Cells.Item(132, 11).Value = 4138168.5
$ws.Cells.Item(132, 13).Value = -4135638.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 6212038
$ws.Cells.Item(20, 9).Value = 7143629
$ws.Cells.Item(20, 11).Value = 7143629
$ws.Cells.Item(20, 13).Value = -7143382

$ws.Cells.Item(86, 8).Value = 4285.5713
$ws.Cells.Item(86, 10).Value = 13600
$ws.Cells.Item(86, 12).Value = 13600
$ws.Cells.Item(86, 14).Value = -15846

$ws.Cells.Item(89, 8).Value = 4285.5713
$ws.Cells.Item(89, 10).Value = 13600
$ws.Cells.Item(89, 12).Value = 68000
$ws.Cells.Item(89, 14).Value = -79232

$ws.Cells.Item(99, 8).Value = 8116.068
$ws.Cells.Item(99, 9).Value = 7990.6665
$ws.Cells.Item(99, 10).Value = 8312.348
$ws.Cells.Item(99, 11).Value = 7990.6665
$ws.Cells.Item(99, 12).Value = 8312.348
$ws.Cells.Item(99, 13).Value = -6492.6665
$ws.Cells.Item(99, 14).Value = -11308.348

$ws.Cells.Item(107, 8).Value = 5556215
$ws.Cells.Item(107, 10).Value = 537.5
$ws.Cells.Item(107, 12).Value = 537.5
$ws.Cells.Item(107, 14).Value = -4377.5

$ws.Cells.Item(134, 8).Value = 1449386
$ws.Cells.Item(134, 9).Value = 1769058.2
$ws.Cells.Item(134, 11).Value = 5307174.6
$ws.Cells.Item(134, 13).Value = -5304639.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 41675596
$ws.Cells.Item(58, 10).Value = 11613.667
$ws.Cells.Item(58, 12).Value = 11613.667
$ws.Cells.Item(58, 14).Value = -12019.667

$ws.Cells.Item(105, 8).Value = 55559676
$ws.Cells.Item(105, 9).Value = 71432970
$ws.Cells.Item(105, 10).Value = 3148.5
$ws.Cells.Item(105, 11).Value = 71432970
$ws.Cells.Item(105, 12).Value = 3148.5
$ws.Cells.Item(105, 13).Value = -71431223
$ws.Cells.Item(105, 14).Value = -6642.5

$ws.Cells.Item(136, 8).Value = 41675596
$ws.Cells.Item(136, 10).Value = 11613.667
$ws.Cells.Item(136, 12).Value = 34841.001
$ws.Cells.Item(136, 14).Value = -39941.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 5000060
$ws.Cells.Item(12, 9).Value = 10000034
$ws.Cells.Item(12, 10).Value = 85
$ws.Cells.Item(12, 11).Value = 30000102
$ws.Cells.Item(12, 12).Value = 255
$ws.Cells.Item(12, 13).Value = -29999929
$ws.Cells.Item(12, 14).Value = -601

$ws.Cells.Item(13, 8).Value = 237.8
$ws.Cells.Item(13, 9).Value = 108.5
$ws.Cells.Item(13, 11).Value = 325.5
$ws.Cells.Item(13, 13).Value = -157.5

$ws.Cells.Item(64, 8).Value = 13616.363
$ws.Cells.Item(64, 10).Value = 15716
$ws.Cells.Item(64, 12).Value = 47148
$ws.Cells.Item(64, 14).Value = -47688

$ws.Cells.Item(67, 8).Value = 13616.363
$ws.Cells.Item(67, 10).Value = 15716
$ws.Cells.Item(67, 12).Value = 47148
$ws.Cells.Item(67, 14).Value = -49020

$ws.Cells.Item(121, 8).Value = 59774.09
$ws.Cells.Item(121, 9).Value = 609.6
$ws.Cells.Item(121, 10).Value = 77175.414
$ws.Cells.Item(121, 11).Value = 1828.8
$ws.Cells.Item(121, 12).Value = 231526.242
$ws.Cells.Item(121, 13).Value = -518.8000000000002
$ws.Cells.Item(121, 14).Value = -234146.242

$ws.Cells.Item(129, 8).Value = 9808532
$ws.Cells.Item(129, 9).Value = 716
$ws.Cells.Item(129, 10).Value = 16674003
$ws.Cells.Item(129, 11).Value = 2148
$ws.Cells.Item(129, 12).Value = 50022009
$ws.Cells.Item(129, 13).Value = 2852
$ws.Cells.Item(129, 14).Value = -50032009

$ws.Cells.Item(132, 8).Value = 2444.5
$ws.Cells.Item(132, 10).Value = 4204
$ws.Cells.Item(132, 12).Value = 37836
$ws.Cells.Item(132, 14).Value = -42896

$ws.Cells.Item(139, 8).Value = 23811736
$ws.Cells.Item(139, 10).Value = 3390.25
$ws.Cells.Item(139, 12).Value = 10170.75
$ws.Cells.Item(139, 14).Value = -20450.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9739.056
$ws.Cells.Item(70, 9).Value = 7390.0835
$ws.Cells.Item(70, 11).Value = 7390.0835
$ws.Cells.Item(70, 13).Value = -7120.0835

$ws.Cells.Item(73, 8).Value = 9739.056
$ws.Cells.Item(73, 9).Value = 7390.0835
$ws.Cells.Item(73, 11).Value = 7390.0835
$ws.Cells.Item(73, 13).Value = -6454.0835

$ws.Cells.Item(113, 8).Value = 8559.076999999999
$ws.Cells.Item(113, 9).Value = 2923.3333
$ws.Cells.Item(113, 11).Value = 2923.3333
$ws.Cells.Item(113, 13).Value = -753.3332999999998

$ws.Cells.Item(123, 8).Value = 44932.332
$ws.Cells.Item(123, 10).Value = 44932.332
$ws.Cells.Item(123, 12).Value = 44932.332
$ws.Cells.Item(123, 14).Value = -49832.332

$ws.Cells.Item(132, 8).Value = 47621340
$ws.Cells.Item(132, 9).Value = 55557610
$ws.Cells.Item(132, 11).Value = 166672830
$ws.Cells.Item(132, 13).Value = -166670300

$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2006.3182
$ws.Cells.Item(16, 9).Value = 690
$ws.Cells.Item(16, 10).Value = 8964
$ws.Cells.Item(16, 11).Value = 690
$ws.Cells.Item(16, 12).Value = 8964
$ws.Cells.Item(16, 13).Value = -520
$ws.Cells.Item(16, 14).Value = -9304

$ws.Cells.Item(40, 8).Value = 5156.4546
$ws.Cells.Item(40, 10).Value = 7624.1
$ws.Cells.Item(40, 12).Value = 7624.1
$ws.Cells.Item(40, 14).Value = -7896.1

$ws.Cells.Item(44, 8).Value = 39998
$ws.Cells.Item(44, 10).Value = 39998
$ws.Cells.Item(44, 12).Value = 39998
$ws.Cells.Item(44, 14).Value = -40910

$ws.Cells.Item(57, 8).Value = 20041
$ws.Cells.Item(57, 9).Value = 20041
$ws.Cells.Item(57, 11).Value = 20041
$ws.Cells.Item(57, 13).Value = -19475

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 14302929
$ws.Cells.Item(13, 10).Value = 24100
$ws.Cells.Item(13, 12).Value = 24100
$ws.Cells.Item(13, 14).Value = -24380

$ws.Cells.Item(46, 8).Value = 85602.164
$ws.Cells.Item(46, 10).Value = 85602.164
$ws.Cells.Item(46, 12).Value = 85602.164
$ws.Cells.Item(46, 14).Value = -86064.164

$ws.Cells.Item(96, 8).Value = 3827.1667
$ws.Cells.Item(96, 9).Value = 2875
$ws.Cells.Item(96, 10).Value = 4303.25
$ws.Cells.Item(96, 11).Value = 2875
$ws.Cells.Item(96, 12).Value = 4303.25
$ws.Cells.Item(96, 13).Value = -1502
$ws.Cells.Item(96, 14).Value = -7049.25

$ws.Cells.Item(107, 8).Value = 554.1818
$ws.Cells.Item(107, 10).Value = 635
$ws.Cells.Item(107, 12).Value = 1905
$ws.Cells.Item(107, 14).Value = -5745

$ws.Cells.Item(113, 8).Value = 6173675
$ws.Cells.Item(113, 10).Value = 461.2
$ws.Cells.Item(113, 12).Value = 1383.6
$ws.Cells.Item(113, 14).Value = -5723.6

$ws.Cells.Item(119, 8).Value = 41289.8
$ws.Cells.Item(119, 10).Value = 41289.8
$ws.Cells.Item(119, 12).Value = 41289.8
$ws.Cells.Item(119, 14).Value = -50965.8

$ws.Cells.Item(122, 8).Value = 1764.826
$ws.Cells.Item(122, 9).Value = 1723.8125
$ws.Cells.Item(122, 11).Value = 5171.4375
$ws.Cells.Item(122, 13).Value = -2721.4375

$ws.Cells.Item(126, 8).Value = 5044.722
$ws.Cells.Item(126, 9).Value = 4196.4
$ws.Cells.Item(126, 11).Value = 12589.2
$ws.Cells.Item(126, 13).Value = -10119.2

$ws.Cells.Item(132, 8).Value = 8434.75
$ws.Cells.Item(132, 9).Value = 6036.2905
$ws.Cells.Item(132, 11).Value = 18108.8715
$ws.Cells.Item(132, 13).Value = -15578.8715

$ws.Cells.Item(134, 8).Value = 85602.164
$ws.Cells.Item(134, 10).Value = 85602.164
$ws.Cells.Item(134, 12).Value = 256806.492
$ws.Cells.Item(134, 14).Value = -261876.492

$ws.Cells.Item(136, 8).Value = 11909559
$ws.Cells.Item(136, 9).Value = 17858026
$ws.Cells.Item(136, 11).Value = 53574078
$ws.Cells.Item(136, 13).Value = -53571528
